$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "57.959.94"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.82%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.969.21"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.59%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "557.26"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.70%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "133.74"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.85%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.515"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.32%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.959.60"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.80%  "
$ws.Range("E10").Value = "  -2.44%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "4.85"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.76%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.450"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.64%  "
$ws.Range("E13").Value = "  +0.59%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.07"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.86%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.456.81"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.60%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.80"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +9.65%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.964.73"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.69%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "58.009.88"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.64%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "418.76"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.82%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.21"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.07%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.686"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.87%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.99"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.26%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.04"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.79%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "79.54"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.71%  "
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.22%  "
$ws.Range("E28").Value = "  -1.87%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.59"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.45%  "
$ws.Range("E30").Value = "  +6.05%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.29"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.01%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.08"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.75%  "
$ws.Range("E33").Value = "  +7.96%  "
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.67"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.56%  "
$ws.Range("B35").Value = "Stacks"
$ws.Range("C35").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.12"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.69%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.940"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.98%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0₃0696"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.68%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "48.56"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.94%  "
$ws.Range("E39").Value = "  +6.22%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.58"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.27%  "
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.109"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.00%  "
$ws.Range("B42").Value = "Bittensor"
$ws.Range("C42").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "382.02"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.48%  "
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0351"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.20%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.681.94"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.92%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.242"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.26%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "122.62"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.55%  "
$ws.Range("E48").Value = "  +2.89%  "
$ws.Range("E49").Value = "  +1.03%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.56"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.26%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.02"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.47%  "
